$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (res) values to "Invalid" for rows 2-6, removing "Valid"
$ws.Range("C2").Value = "Invalid"
$ws.Range("C3").Value = "Invalid"
$ws.Range("C4").Value = "Invalid"
$ws.Range("C5").Value = "Invalid"
$ws.Range("C6").Value = "Invalid"

# Update the selected cell to C3
$ws.Range("C3").Select()
